$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("B2:G21")
$rng.NumberFormat = "@"

$ws.Cells.Item(2, 2).Value = "Galatasaray"
$ws.Cells.Item(2, 3).Value = "1.8"
$ws.Cells.Item(2, 4).Value = "7.0"
$ws.Cells.Item(2, 5).Value = "78%"
$ws.Cells.Item(2, 6).Value = "60%"
$ws.Cells.Item(2, 7).Value = "2.97"
$ws.Cells.Item(3, 2).Value = "Fenerbahçe"
$ws.Cells.Item(3, 3).Value = "2.5"
$ws.Cells.Item(3, 4).Value = "6.0"
$ws.Cells.Item(3, 5).Value = "82%"
$ws.Cells.Item(3, 6).Value = "67%"
$ws.Cells.Item(3, 7).Value = "3.55"
$ws.Cells.Item(4, 2).Value = "Trabzonspor"
$ws.Cells.Item(4, 3).Value = "1.8"
$ws.Cells.Item(4, 4).Value = "5.2"
$ws.Cells.Item(4, 5).Value = "72%"
$ws.Cells.Item(4, 6).Value = "53%"
$ws.Cells.Item(4, 7).Value = "3.00"
$ws.Cells.Item(5, 2).Value = "Beşiktaş"
$ws.Cells.Item(5, 3).Value = "2.3"
$ws.Cells.Item(5, 4).Value = "5.5"
$ws.Cells.Item(5, 5).Value = "82%"
$ws.Cells.Item(5, 6).Value = "43%"
$ws.Cells.Item(5, 7).Value = "2.48"
$ws.Cells.Item(6, 2).Value = "Başakşehir"
$ws.Cells.Item(6, 3).Value = "2.2"
$ws.Cells.Item(6, 4).Value = "4.2"
$ws.Cells.Item(6, 5).Value = "70%"
$ws.Cells.Item(6, 6).Value = "36%"
$ws.Cells.Item(6, 7).Value = "2.48"
$ws.Cells.Item(7, 2).Value = "Rizespor"
$ws.Cells.Item(7, 3).Value = "2.3"
$ws.Cells.Item(7, 4).Value = "4.7"
$ws.Cells.Item(7, 5).Value = "70%"
$ws.Cells.Item(7, 6).Value = "55%"
$ws.Cells.Item(7, 7).Value = "2.79"
$ws.Cells.Item(8, 2).Value = "Kasımpaşa"
$ws.Cells.Item(8, 3).Value = "1.9"
$ws.Cells.Item(8, 4).Value = "5.1"
$ws.Cells.Item(8, 5).Value = "88%"
$ws.Cells.Item(8, 6).Value = "67%"
$ws.Cells.Item(8, 7).Value = "3.45"
$ws.Cells.Item(9, 2).Value = "Antalyaspor"
$ws.Cells.Item(9, 3).Value = "2.1"
$ws.Cells.Item(9, 4).Value = "5.5"
$ws.Cells.Item(9, 5).Value = "73%"
$ws.Cells.Item(9, 6).Value = "40%"
$ws.Cells.Item(9, 7).Value = "2.30"
$ws.Cells.Item(10, 2).Value = "Alanyaspor"
$ws.Cells.Item(10, 3).Value = "2.5"
$ws.Cells.Item(10, 4).Value = "4.5"
$ws.Cells.Item(10, 5).Value = "76%"
$ws.Cells.Item(10, 6).Value = "55%"
$ws.Cells.Item(10, 7).Value = "2.70"
$ws.Cells.Item(11, 2).Value = "Sivasspor"
$ws.Cells.Item(11, 3).Value = "2.1"
$ws.Cells.Item(11, 4).Value = "3.8"
$ws.Cells.Item(11, 5).Value = "69%"
$ws.Cells.Item(11, 6).Value = "45%"
$ws.Cells.Item(11, 7).Value = "2.58"
$ws.Cells.Item(12, 2).Value = "Adana Demirspor"
$ws.Cells.Item(12, 3).Value = "2.2"
$ws.Cells.Item(12, 4).Value = "5.0"
$ws.Cells.Item(12, 5).Value = "73%"
$ws.Cells.Item(12, 6).Value = "55%"
$ws.Cells.Item(12, 7).Value = "2.82"
$ws.Cells.Item(13, 2).Value = "Samsunspor"
$ws.Cells.Item(13, 3).Value = "2.0"
$ws.Cells.Item(13, 4).Value = "4.8"
$ws.Cells.Item(13, 5).Value = "82%"
$ws.Cells.Item(13, 6).Value = "39%"
$ws.Cells.Item(13, 7).Value = "2.42"
$ws.Cells.Item(14, 2).Value = "Ankaragücü"
$ws.Cells.Item(14, 3).Value = "1.5"
$ws.Cells.Item(14, 4).Value = "3.7"
$ws.Cells.Item(14, 5).Value = "82%"
$ws.Cells.Item(14, 6).Value = "43%"
$ws.Cells.Item(14, 7).Value = "2.52"
$ws.Cells.Item(15, 2).Value = "Kayserispor"
$ws.Cells.Item(15, 3).Value = "2.0"
$ws.Cells.Item(15, 4).Value = "4.4"
$ws.Cells.Item(15, 5).Value = "76%"
$ws.Cells.Item(15, 6).Value = "43%"
$ws.Cells.Item(15, 7).Value = "2.55"
$ws.Cells.Item(16, 2).Value = "Konyaspor"
$ws.Cells.Item(16, 3).Value = "1.9"
$ws.Cells.Item(16, 4).Value = "4.1"
$ws.Cells.Item(16, 5).Value = "85%"
$ws.Cells.Item(16, 6).Value = "33%"
$ws.Cells.Item(16, 7).Value = "2.45"
$ws.Cells.Item(17, 2).Value = "Gaziantep"
$ws.Cells.Item(17, 3).Value = "2.2"
$ws.Cells.Item(17, 4).Value = "4.2"
$ws.Cells.Item(17, 5).Value = "85%"
$ws.Cells.Item(17, 6).Value = "49%"
$ws.Cells.Item(17, 7).Value = "2.64"
$ws.Cells.Item(18, 2).Value = "Karagümrük"
$ws.Cells.Item(18, 3).Value = "2.6"
$ws.Cells.Item(18, 4).Value = "5.1"
$ws.Cells.Item(18, 5).Value = "76%"
$ws.Cells.Item(18, 6).Value = "42%"
$ws.Cells.Item(18, 7).Value = "2.39"
$ws.Cells.Item(19, 2).Value = "Hatayspor"
$ws.Cells.Item(19, 3).Value = "2.8"
$ws.Cells.Item(19, 4).Value = "4.4"
$ws.Cells.Item(19, 5).Value = "73%"
$ws.Cells.Item(19, 6).Value = "49%"
$ws.Cells.Item(19, 7).Value = "2.52"
$ws.Cells.Item(20, 2).Value = "Pendikspor"
$ws.Cells.Item(20, 3).Value = "2.2"
$ws.Cells.Item(20, 4).Value = "4.7"
$ws.Cells.Item(20, 5).Value = "88%"
$ws.Cells.Item(20, 6).Value = "51%"
$ws.Cells.Item(20, 7).Value = "3.18"
$ws.Cells.Item(21, 2).Value = "İstanbulspor"
$ws.Cells.Item(21, 3).Value = "2.5"
$ws.Cells.Item(21, 4).Value = "4.0"
$ws.Cells.Item(21, 5).Value = "78%"
$ws.Cells.Item(21, 6).Value = "50%"
$ws.Cells.Item(21, 7).Value = "2.63"

$rng.Style = "Normal"